$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 257
$ws.Range("F5").Value = 288
$ws.Range("F7").Value = 365
$ws.Range("F8").Value = 1854
$ws.Range("F9").Value = 799
$ws.Range("F10").Value = 17
$ws.Range("F11").Value = 19
$ws.Range("F12").Value = 1561
$ws.Range("F13").Value = 1561
$ws.Range("F14").Value = 1299
$ws.Range("F15").Value = 34
$ws.Range("F16").Value = 1371
$ws.Range("F17").Value = 174
$ws.Range("F18").Value = 386
$ws.Range("F21").Value = 125
$ws.Range("F22").Value = 6823
$ws.Range("F23").Value = 7300
$ws.Range("F24").Value = 17
$ws.Range("F25").Value = 163
$ws.Range("F26").Value = 489
$ws.Range("F28").Value = 225
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = 4
$ws.Range("F33").Value = 1341
$ws.Range("F36").Value = 641
$ws.Range("F39").Value = 289
$ws.Range("F40").Value = 161
$ws.Range("F41").Value = 166
$ws.Range("F42").Value = 73
$ws.Range("F44").Value = 112

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 16
$ws.Range("F9").Value = 18
$ws.Range("F17").Value = 260

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 233
$ws.Range("F5").Value = 97

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 257
$ws.Range("F6").Value = 233
$ws.Range("F7").Value = 97
$ws.Range("F9").Value = 288
$ws.Range("F10").Value = 16
$ws.Range("F12").Value = 365
$ws.Range("F13").Value = 1854
$ws.Range("F14").Value = 799
$ws.Range("F15").Value = 17
$ws.Range("F16").Value = 19
$ws.Range("F17").Value = 1561
$ws.Range("F18").Value = 1561
$ws.Range("F19").Value = 1299
$ws.Range("F20").Value = 34
$ws.Range("F21").Value = 1371
$ws.Range("F22").Value = 174
$ws.Range("F23").Value = 386
$ws.Range("F25").Value = 125
$ws.Range("F27").Value = 6823
$ws.Range("F28").Value = 7300
$ws.Range("F29").Value = 163
$ws.Range("F30").Value = 225
$ws.Range("F31").Value = 1341
$ws.Range("F35").Value = 18
$ws.Range("F38").Value = 641
$ws.Range("F43").Value = 289
$ws.Range("F44").Value = 166
$ws.Range("F45").Value = 73
$ws.Range("F46").Value = 96
$ws.Range("F47").Value = 112
$ws.Range("F49").Value = 260
